$d = $word.ActiveDocument

# --- Change 1 (Layout paragraph): tighten the two-column sentence ---
$d.Content.Find.Execute(
    "in a two-column, full-justified format. They break the text into two columns, which reduces",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "in a two-column format. This reduces", 2) | Out-Null

# --- Change 2 (Typography paragraph): drop the body-text-typeface sentence ---
$d.Content.Find.Execute(
    "the selected page. They set the body text in an upright serif face. The typesetters italicize",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "the selected page. The typesetters italicize", 2) | Out-Null

# --- Change 3 (Line Drawings paragraph): merge the figure-number sentences ---
# The existing "_GoBack" bookmark sits between "number," and " which the body text
# uses...". A Find/Replace that spans that boundary drops the bookmark outright, so
# rewrite the wording first (which removes the stale bookmark along with the text
# around it), then re-add "_GoBack" collapsed at the very end of the paragraph once
# the final wording is in place.
$d.Content.Find.Execute(
    "The typesetters also include a figure number, which the body text uses to refer to the visual aid.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The typesetters include a figure number, which allows the body text to incorporate the drawing while preserving the flow of prose.", 2) | Out-Null

# Re-anchor "_GoBack" at the end of the rewritten sentence. A bookmark collapsed
# exactly at a paragraph's final character position lands in the wrong spot in this
# runtime, so temporarily append a marker character, bookmark just before it, then
# remove the marker again.
$d.Content.Find.Execute(
    "preserving the flow of prose.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "preserving the flow of prose.#", 2) | Out-Null

$full = $d.Content.Text
$markerPos = $full.IndexOf("preserving the flow of prose.#") + ("preserving the flow of prose.").Length
$anchor = $d.Range($markerPos, $markerPos)
$d.Bookmarks.Add("_GoBack", $anchor) | Out-Null

$d.Range($markerPos, $markerPos + 1).Delete() | Out-Null
